$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 650
$ws.Range("I29").Value = 650
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1950
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -1669
$ws.Range("N29").ClearContents()
$ws.Range("H86").Value = 121215880
$ws.Range("I86").Value = 100004130
$ws.Range("K86").Value = 100004130
$ws.Range("M86").Value = -100003007
$ws.Range("H89").Value = 121215880
$ws.Range("I89").Value = 100004130
$ws.Range("K89").Value = 500020650
$ws.Range("M89").Value = -500015034
$ws.Range("H116").Value = 11021.174
$ws.Range("I116").Value = 15780.571
$ws.Range("J116").Value = 3617.6667
$ws.Range("K116").Value = 15780.571
$ws.Range("L116").Value = 3617.6667
$ws.Range("M116").Value = -12338.571
$ws.Range("N116").Value = -10501.6667
$ws.Range("H118").Value = 808.6667
$ws.Range("I118").Value = 677.5833
$ws.Range("K118").Value = 2032.7499
$ws.Range("M118").Value = -375.7499
$ws.Range("H137").Value = 9092933
$ws.Range("I137").Value = 2225.8
$ws.Range("K137").Value = 6677.400000000001
$ws.Range("M137").Value = -4127.400000000001
$ws.Range("H138").Value = 11666.26
$ws.Range("I138").Value = 14587.5
$ws.Range("J138").Value = 11544.542
$ws.Range("K138").Value = 43762.5
$ws.Range("L138").Value = 34633.626
$ws.Range("M138").Value = -38622.5
$ws.Range("N138").Value = -44913.626

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4598074
$ws.Range("I32").Value = 5928944
$ws.Range("K32").Value = 5928944
$ws.Range("M32").Value = -5928657
$ws.Range("H61").Value = 2225504
$ws.Range("I61").Value = 12271.5
$ws.Range("J61").Value = 4608985.5
$ws.Range("K61").Value = 12271.5
$ws.Range("L61").Value = 4608985.5
$ws.Range("M61").Value = -12059.5
$ws.Range("N61").Value = -4609409.5
$ws.Range("H74").Value = 327116.97
$ws.Range("I74").Value = 4351.148
$ws.Range("K74").Value = 4351.148
$ws.Range("M74").Value = -3477.148
$ws.Range("H77").Value = 327116.97
$ws.Range("I77").Value = 4351.148
$ws.Range("K77").Value = 21755.74
$ws.Range("M77").Value = -17387.74
$ws.Range("H97").Value = 8628.4
$ws.Range("I97").Value = 8887.571
$ws.Range("K97").Value = 8887.571
$ws.Range("M97").Value = -8391.571
$ws.Range("H136").Value = 2225504
$ws.Range("I136").Value = 12271.5
$ws.Range("J136").Value = 4608985.5
$ws.Range("K136").Value = 36814.5
$ws.Range("L136").Value = 13826956.5
$ws.Range("M136").Value = -34264.5
$ws.Range("N136").Value = -13832056.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1789.0344
$ws.Range("I20").Value = 1624.4667
$ws.Range("J20").Value = 1965.3572
$ws.Range("K20").Value = 1624.4667
$ws.Range("L20").Value = 1965.3572
$ws.Range("M20").Value = -1377.4667
$ws.Range("N20").Value = -2459.3572
$ws.Range("H86").Value = 3300.5
$ws.Range("I86").Value = 2117.3333
$ws.Range("J86").Value = 4187.875
$ws.Range("K86").Value = 2117.3333
$ws.Range("L86").Value = 4187.875
$ws.Range("M86").Value = -994.3332999999998
$ws.Range("N86").Value = -6433.875
$ws.Range("H89").Value = 3300.5
$ws.Range("I89").Value = 2117.3333
$ws.Range("J89").Value = 4187.875
$ws.Range("K89").Value = 10586.6665
$ws.Range("L89").Value = 20939.375
$ws.Range("M89").Value = -4970.666499999999
$ws.Range("N89").Value = -32171.375
$ws.Range("H105").Value = 11912582
$ws.Range("I105").Value = 18526430
$ws.Range("J105").Value = 7655.7
$ws.Range("K105").Value = 18526430
$ws.Range("L105").Value = 7655.7
$ws.Range("M105").Value = -18524683
$ws.Range("N105").Value = -11149.7
$ws.Range("H134").Value = 45054380
$ws.Range("I134").Value = 66850.69
$ws.Range("K134").Value = 200552.07
$ws.Range("M134").Value = -198017.07

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2744.389
$ws.Range("I22").Value = 552.8823
$ws.Range("K22").Value = 552.8823
$ws.Range("M22").Value = -202.8823
$ws.Range("H31").Value = 4471.048
$ws.Range("I31").Value = 3925.6
$ws.Range("K31").Value = 3925.6
$ws.Range("M31").Value = -3630.6
$ws.Range("H34").Value = 4471.048
$ws.Range("I34").Value = 3925.6
$ws.Range("K34").Value = 3925.6
$ws.Range("M34").Value = -3723.6
$ws.Range("H58").Value = 4148.15
$ws.Range("I58").Value = 1810.3334
$ws.Range("J58").Value = 5150.0713
$ws.Range("K58").Value = 1810.3334
$ws.Range("L58").Value = 5150.0713
$ws.Range("M58").Value = -1607.3334
$ws.Range("N58").Value = -5556.0713
$ws.Range("H86").Value = 13064.4
$ws.Range("I86").Value = 6439.25
$ws.Range("J86").Value = 15473.546
$ws.Range("K86").Value = 6439.25
$ws.Range("L86").Value = 15473.546
$ws.Range("M86").Value = -5316.25
$ws.Range("N86").Value = -17719.546
$ws.Range("H89").Value = 13064.4
$ws.Range("I89").Value = 6439.25
$ws.Range("J89").Value = 15473.546
$ws.Range("K89").Value = 32196.25
$ws.Range("L89").Value = 77367.73
$ws.Range("M89").Value = -26580.25
$ws.Range("N89").Value = -88599.73
$ws.Range("H136").Value = 4148.15
$ws.Range("I136").Value = 1810.3334
$ws.Range("J136").Value = 5150.0713
$ws.Range("K136").Value = 5431.0002
$ws.Range("L136").Value = 15450.2139
$ws.Range("M136").Value = -2881.0002
$ws.Range("N136").Value = -20550.2139
$ws.Range("H141").Value = 559545.5
$ws.Range("J141").Value = 588120.1
$ws.Range("L141").Value = 588120.1
$ws.Range("N141").Value = -598480.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 113782.89
$ws.Range("I33").Value = 350
$ws.Range("J33").Value = 170499.33
$ws.Range("K33").Value = 2100
$ws.Range("L33").Value = 1022995.98
$ws.Range("M33").Value = -1817
$ws.Range("N33").Value = -1023561.98
$ws.Range("H81").Value = 7418
$ws.Range("I81").Value = 4276.4
$ws.Range("K81").Value = 12829.2
$ws.Range("M81").Value = -11706.2
$ws.Range("H84").Value = 7418
$ws.Range("I84").Value = 4276.4
$ws.Range("K84").Value = 38487.6
$ws.Range("M84").Value = -32871.6
$ws.Range("H107").Value = 715.1111
$ws.Range("I107").Value = 558.0909
$ws.Range("J107").Value = 823.0625
$ws.Range("K107").Value = 1674.2727
$ws.Range("L107").Value = 2469.1875
$ws.Range("M107").Value = 245.7273
$ws.Range("N107").Value = -6309.1875
$ws.Range("H113").Value = 2257.9443
$ws.Range("I113").Value = 2184.75
$ws.Range("J113").Value = 2278.8572
$ws.Range("K113").Value = 6554.25
$ws.Range("L113").Value = 6836.571599999999
$ws.Range("M113").Value = -4384.25
$ws.Range("N113").Value = -11176.5716
$ws.Range("H128").Value = 239998
$ws.Range("I128").Value = 239998
$ws.Range("K128").Value = 719994
$ws.Range("M128").Value = -715014

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 45737680
$ws.Range("I80").Value = 269077.84
$ws.Range("J80").Value = 100300000
$ws.Range("K80").Value = 269077.84
$ws.Range("L80").Value = 100300000
$ws.Range("M80").Value = -268079.84
$ws.Range("N80").Value = -100301996
$ws.Range("H83").Value = 45737680
$ws.Range("I83").Value = 269077.84
$ws.Range("J83").Value = 100300000
$ws.Range("K83").Value = 1345389.2
$ws.Range("L83").Value = 501500000
$ws.Range("M83").Value = -1340397.2
$ws.Range("N83").Value = -501509984
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H97").Value = 75128.96000000001
$ws.Range("I97").Value = 42803.543
$ws.Range("J97").Value = 333732.34
$ws.Range("K97").Value = 42803.543
$ws.Range("L97").Value = 333732.34
$ws.Range("M97").Value = -42307.543
$ws.Range("N97").Value = -334724.34

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7424.4546
$ws.Range("I136").Value = 7726.077
$ws.Range("K136").Value = 23178.231
$ws.Range("M136").Value = -20628.231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11729.9
$ws.Range("J81").Value = 2400
$ws.Range("L81").Value = 4800
$ws.Range("N81").Value = -6922
$ws.Range("H84").Value = 11729.9
$ws.Range("J84").Value = 2400
$ws.Range("L84").Value = 24000
$ws.Range("N84").Value = -34608
$ws.Range("H113").Value = 166667140
$ws.Range("I113").Value = 250000370
$ws.Range("J113").Value = 675
$ws.Range("K113").Value = 750001110
$ws.Range("L113").Value = 2025
$ws.Range("M113").Value = -749998940
$ws.Range("N113").Value = -6365
$ws.Range("H132").Value = 41660.6
$ws.Range("I132").Value = 134769.14
$ws.Range("J132").Value = 1756.9429
$ws.Range("K132").Value = 404307.42
$ws.Range("L132").Value = 5270.8287
$ws.Range("M132").Value = -401777.42
$ws.Range("N132").Value = -10330.8287
